$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: _old -> _FV2210, _new -> _FV2304 --------------
$oldNames = @(
    "Segmentname_old", "Segmentgruppe_old", "Segment_old", "Datenelement_old",
    "Segment ID_old", "Code_old", "Qualifier_old", "Beschreibung_old",
    "Bedingungsausdruck_old", "Bedingung_old"
)
$newNames = @(
    "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210",
    "Segment ID_FV2210", "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210", "Bedingung_FV2210"
)
for ($i = 0; $i -lt $oldNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newNames[$i]
}

$newOldNames = @(
    "Segmentname_new", "Segmentgruppe_new", "Segment_new", "Datenelement_new",
    "Segment ID_new", "Code_new", "Qualifier_new", "Beschreibung_new",
    "Bedingungsausdruck_new", "Bedingung_new"
)
$newNewNames = @(
    "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304",
    "Segment ID_FV2304", "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)
for ($i = 0; $i -lt $newOldNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newNewNames[$i]
}

# --- Add table over the used data range (A1:U80) -------------------------
$tableRange = $ws.Range("A1:U80")
$listObj = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObj.Name = "Table1"

# --- Freeze header row -----------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
